$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8559919436052367
$ws.Range("B4").Value = 0.8619201725997843
$ws.Range("B5").Value = 0.2849162011173184
$ws.Range("B6").Value = 0.9815724815724816
$ws.Range("B7").Value = 0.8968324016867707
